$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.309.31"
$ws.Range("E2").Value = "  -0.42%  "

$ws.Range("D3").Value = "2.062.98"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("D5").Value = "'233.71"
$ws.Range("E5").Value = "  -0.84%  "

$ws.Range("D6").Value = "'0.622"
$ws.Range("E6").Value = "  +1.08%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'56.86"
$ws.Range("E8").Value = "  -1.57%  "

$ws.Range("D9").Value = "'0.386"
$ws.Range("E9").Value = "  +1.10%  "

$ws.Range("E10").Value = "  +0.44%  "

$ws.Range("E11").Value = "  +0.81%  "

$ws.Range("D12").Value = "2.367.01"
$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("D13").Value = "'14.41"
$ws.Range("E13").Value = "  +1.16%  "

$ws.Range("D14").Value = "'20.66"
$ws.Range("E14").Value = "  -0.78%  "

$ws.Range("E15").Value = "  -0.25%  "

$ws.Range("E16").Value = "  -0.50%  "

$ws.Range("D17").Value = "2.062.71"
$ws.Range("E17").Value = "  -0.30%  "

$ws.Range("D18").Value = "37.231.44"
$ws.Range("E18").Value = "  -0.87%  "

$ws.Range("D19").Value = "'6.36"
$ws.Range("E19").Value = "  +3.99%  "

$ws.Range("D20").Value = "'69.43"
$ws.Range("E20").Value = "  +1.22%  "

$ws.Range("D21").Value = "0.0₃0813"
$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("D22").Value = "'226.21"
$ws.Range("E22").Value = "  +0.47%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("E24").Value = "  +0.93%  "

$ws.Range("D25").Value = "'2.39"
$ws.Range("E25").Value = "  -2.54%  "

$ws.Range("D26").Value = "'166.16"
$ws.Range("E26").Value = "  +1.62%  "

$ws.Range("E27").Value = "  -0.84%  "

$ws.Range("E28").Value = "  +2.21%  "

$ws.Range("D29").Value = "'18.98"
$ws.Range("E29").Value = "  -1.01%  "

$ws.Range("E30").Value = "  -2.38%  "

$ws.Range("D31").Value = "'0.117"
$ws.Range("E31").Value = "  -0.93%  "

$ws.Range("E32").Value = "  +0.21%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0616"
$ws.Range("E33").Value = "  -1.58%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.60"
$ws.Range("E34").Value = "  +3.45%  "

$ws.Range("E35").Value = "  -4.73%  "

$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("E37").Value = "  -0.27%  "

$ws.Range("D38").Value = "'3.21"
$ws.Range("E38").Value = "  -4.36%  "

$ws.Range("E39").Value = "  -5.21%  "

$ws.Range("D40").Value = "'2.95"
$ws.Range("E40").Value = "  -0.75%  "

$ws.Range("E41").Value = "  +0.34%  "

$ws.Range("D42").Value = "'95.98"
$ws.Range("E42").Value = "  +0.62%  "

$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").Value = "'4.31"
$ws.Range("E44").Value = "  -3.65%  "

$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'1.17"
$ws.Range("E45").Value = "  +3.24%  "

$ws.Range("D46").Value = "'0.0927"
$ws.Range("E46").Value = "  -3.06%  "

$ws.Range("E47").Value = "  -0.32%  "

$ws.Range("D48").Value = "'15.13"
$ws.Range("E48").Value = "  -5.98%  "

$ws.Range("E49").Value = "  -1.63%  "

$ws.Range("D50").Value = "'2.96"
$ws.Range("E50").Value = "  +0.60%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.256.70"
$ws.Range("E51").Value = "  -0.09%  "
